# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-17 13:30:08
#
# Applies the attendance-sync edit to the "Session Analysis Results" sheet:
#  1. Class Statistics summary (L6, L7, L9, L10) refreshed with new totals.
#  2. Group Statistics table (rows 16-26, cols O/P/R/S) refreshed for the
#     groups whose General Surgery session 5 just got recorded.
#  3. "Recorded By" cells that used to read "<email>, System" now read
#     "System, <email>" (System ran first).
#  4. Six General Surgery "session 5" rows (B1-10, B1-11, B1-12, B1-7, B1-8,
#     B1-9) flip from Not-Recorded/pending to Recorded, with a recorder
#     e-mail, a new attendance fraction, and the "Recorded" (green) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteAll / xlPasteFormats constants (avoid relying on named enum support)
$xlPasteFormats = -4122

function Set-TextCell {
    # Writes a literal text string into a cell without Excel's automatic
    # "looks like a number/percentage" coercion, and without leaving the
    # cell's number format changed (restores it from a same-style sibling
    # cell so the xf stays "General").
    param(
        [object]$Sheet,
        [string]$Addr,
        [string]$Text,
        [string]$FormatDonorAddr
    )
    $Sheet.Range($Addr).NumberFormat = "@"
    $Sheet.Range($Addr).Value2 = $Text
    $Sheet.Range($FormatDonorAddr).Copy()
    $Sheet.Range($Addr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# 1. Class Statistics (K2:L10 block)
# ---------------------------------------------------------------------
$ws.Range("L6").Value2 = 73
$ws.Range("L7").Value2 = 5
Set-TextCell $ws "L9"  "28.3%" "L4"
Set-TextCell $ws "L10" "82.0%" "L4"

# ---------------------------------------------------------------------
# 2. Group Statistics table updates (Recorded / Missing / Coverage % /
#    Avg Attendance % for the groups affected by the new session-5 data)
# ---------------------------------------------------------------------
$groupRows = @{
    16 = @{ O = 6; P = 1; R = "27.3%"; S = "79.6%" }   # B1-10
    17 = @{ O = 6; P = 1; R = "27.3%"; S = "70.2%" }   # B1-11
    18 = @{ O = 7; P = 0; R = "31.8%"; S = "77.6%" }   # B1-12
    24 = @{ O = 6; P = 1; R = "27.3%"; S = "73.5%" }   # B1-7
    25 = @{ O = 6; P = 1; R = "27.3%"; S = "80.5%" }   # B1-8
    26 = @{ O = 6; P = 1; R = "27.3%"; S = "74.7%" }   # B1-9
}
foreach ($row in $groupRows.Keys) {
    $vals = $groupRows[$row]
    $ws.Range("O$row").Value2 = $vals.O
    $ws.Range("P$row").Value2 = $vals.P
    Set-TextCell $ws "R$row" $vals.R "O$row"
    Set-TextCell $ws "S$row" $vals.S "O$row"
}

# ---------------------------------------------------------------------
# 3. "Recorded By" text: "<email>, System" -> "System, <email>"
# ---------------------------------------------------------------------
$recordedByRows = @(2,3,4,23,24,25,26,45,46,47,48,67,68,69,70,89,90,91,110,111,112,131,132,133,152,153,154,173,174,175,194,195,196,197,216,217,218,219,238,239,240,241)
foreach ($row in $recordedByRows) {
    $ws.Range("G$row").Value2 = "System, dnasr281@gmail.com"
}

# ---------------------------------------------------------------------
# 4. General Surgery "session 5" rows that just got recorded.
#    Each becomes Recorder/Students/Status + flips the row look from the
#    "Not Recorded" (pink) style to the "Recorded" (green) style used by
#    row 2 (A2:I2).
# ---------------------------------------------------------------------
$recordedSessionRows = @{
    27  = @{ G = "dnasr281@gmail.com"; H = "24/31"; I = "Recorded" }   # B1-10
    49  = @{ G = "dnasr281@gmail.com"; H = "12/19"; I = "Recorded" }   # B1-11
    71  = @{ G = "dnasr281@gmail.com"; H = "20/21"; I = "Recorded" }   # B1-12
    198 = @{ G = "dnasr281@gmail.com"; H = "22/27"; I = "Recorded" }   # B1-7
    220 = @{ G = "dnasr281@gmail.com"; H = "24/29"; I = "Recorded" }   # B1-8
    242 = @{ G = "dnasr281@gmail.com"; H = "22/29"; I = "Recorded" }   # B1-9
}
foreach ($row in $recordedSessionRows.Keys) {
    $vals = $recordedSessionRows[$row]

    # Repaint A:I with the "Recorded" row look (copied from the A2:I2
    # template row) before writing the new values.
    $ws.Range("A2:I2").Copy()
    $ws.Range("A$row`:I$row").PasteSpecial($xlPasteFormats)

    $ws.Range("G$row").Value2 = $vals.G
    $ws.Range("H$row").Value2 = $vals.H
    $ws.Range("I$row").Value2 = $vals.I
}

$excel.CutCopyMode = $false
